# Update Grade column (B) from "C" to "B" for the students whose grade
# was re-assessed: rows 6, 18, 21, 25, 36, 45 (No. 5, 17, 20, 24, 35, 44)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToUpdate = @(6, 18, 21, 25, 36, 45)
foreach ($r in $rowsToUpdate) {
    $ws.Range("B$r").Value = "B"
}

# Match the view/selection state left behind by the author's edit session
$ws.Range("G16").Select()
